# Update market-data cells (currentAveragePrice / Leve price / profit columns)
# per the scheduled-runner refresh. Values sourced from the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3672.6667
$ws.Range("I74").Value = 3408.9092
$ws.Range("J74").Value = 3804.5454
$ws.Range("K74").Value = 3408.9092
$ws.Range("L74").Value = 3804.5454
$ws.Range("M74").Value = -2472.9092
$ws.Range("N74").Value = -5676.5454
$ws.Range("H77").Value = 3672.6667
$ws.Range("I77").Value = 3408.9092
$ws.Range("J77").Value = 3804.5454
$ws.Range("K77").Value = 17044.546
$ws.Range("L77").Value = 19022.727
$ws.Range("M77").Value = -12364.546
$ws.Range("N77").Value = -28382.727
$ws.Range("H100").Value = 2912.2222
$ws.Range("I100").Value = 700
$ws.Range("J100").Value = 3544.2856
$ws.Range("K100").Value = 700
$ws.Range("L100").Value = 3544.2856
$ws.Range("M100").Value = -159
$ws.Range("N100").Value = -4626.2856
$ws.Range("H125").Value = 4084.8333
$ws.Range("I125").Value = 5402.25
$ws.Range("K125").Value = 48620.25
$ws.Range("M125").Value = -46160.25
$ws.Range("H132").Value = 3635.3333
$ws.Range("I132").Value = 3498.7407
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 10496.2221
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -7966.222099999999
$ws.Range("N132").Value = -17810
$ws.Range("H141").Value = 4290.364
$ws.Range("I141").Value = 2761.3333
$ws.Range("J141").Value = 4863.75
$ws.Range("K141").Value = 8283.999899999999
$ws.Range("L141").Value = 14591.25
$ws.Range("M141").Value = -3103.999899999999
$ws.Range("N141").Value = -24951.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10657.917
$ws.Range("I32").Value = 11190.78
$ws.Range("K32").Value = 11190.78
$ws.Range("M32").Value = -10903.78
$ws.Range("H61").Value = 20002238
$ws.Range("I61").Value = 22729500
$ws.Range("J61").Value = 2304.6667
$ws.Range("K61").Value = 22729500
$ws.Range("L61").Value = 2304.6667
$ws.Range("M61").Value = -22729288
$ws.Range("N61").Value = -2728.6667
$ws.Range("H88").Value = 2400.5908
$ws.Range("J88").Value = 2362.25
$ws.Range("L88").Value = 2362.25
$ws.Range("N88").Value = -3174.25
$ws.Range("H91").Value = 2400.5908
$ws.Range("J91").Value = 2362.25
$ws.Range("L91").Value = 2362.25
$ws.Range("N91").Value = -5170.25
$ws.Range("H102").Value = 2652.5
$ws.Range("I102").Value = 2010
$ws.Range("J102").Value = 2866.6667
$ws.Range("K102").Value = 2010
$ws.Range("L102").Value = 2866.6667
$ws.Range("M102").Value = -388
$ws.Range("N102").Value = -6110.6667
$ws.Range("H128").Value = 47271.75
$ws.Range("J128").Value = 47271.75
$ws.Range("L128").Value = 47271.75
$ws.Range("N128").Value = -57231.75
$ws.Range("H136").Value = 20002238
$ws.Range("I136").Value = 22729500
$ws.Range("J136").Value = 2304.6667
$ws.Range("K136").Value = 68188500
$ws.Range("L136").Value = 6914.000100000001
$ws.Range("M136").Value = -68185950
$ws.Range("N136").Value = -12014.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23811700
$ws.Range("I86").Value = 1767.6666
$ws.Range("K86").Value = 1767.6666
$ws.Range("M86").Value = -644.6666
$ws.Range("H89").Value = 23811700
$ws.Range("I89").Value = 1767.6666
$ws.Range("K89").Value = 1767.6666
$ws.Range("M89").Value = -3222.333000000001
$ws.Range("H99").Value = 938.46155
$ws.Range("I99").Value = 892.8570999999999
$ws.Range("J99").Value = 991.6667
$ws.Range("K99").Value = 892.8570999999999
$ws.Range("L99").Value = 991.6667
$ws.Range("M99").Value = 605.1429000000001
$ws.Range("N99").Value = -3987.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15153037
$ws.Range("I31").Value = 1593.8572
$ws.Range("J31").Value = 333333340
$ws.Range("K31").Value = 1593.8572
$ws.Range("L31").Value = 333333340
$ws.Range("M31").Value = -1298.8572
$ws.Range("N31").Value = -333333930
$ws.Range("H34").Value = 15153037
$ws.Range("I34").Value = 1593.8572
$ws.Range("J34").Value = 333333340
$ws.Range("K34").Value = 1593.8572
$ws.Range("L34").Value = 333333340
$ws.Range("M34").Value = -1391.8572
$ws.Range("N34").Value = -333333744
$ws.Range("H62").Value = 2533.3333
$ws.Range("I62").Value = 2225
$ws.Range("J62").Value = 3150
$ws.Range("K62").Value = 2225
$ws.Range("L62").Value = 3150
$ws.Range("M62").Value = -1601
$ws.Range("N62").Value = -4398
$ws.Range("H65").Value = 2533.3333
$ws.Range("I65").Value = 2225
$ws.Range("J65").Value = 3150
$ws.Range("K65").Value = 11125
$ws.Range("L65").Value = 15750
$ws.Range("M65").Value = -8005
$ws.Range("N65").Value = -21990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1517.3334
$ws.Range("I5").Value = 768
$ws.Range("J5").Value = 2266.6667
$ws.Range("K5").Value = 2304
$ws.Range("L5").Value = 6800.000100000001
$ws.Range("M5").Value = -2192
$ws.Range("N5").Value = -7024.000100000001
$ws.Range("H80").Value = 1845.909
$ws.Range("J80").Value = 1811.2222
$ws.Range("L80").Value = 5433.6666
$ws.Range("N80").Value = -7305.6666
$ws.Range("H83").Value = 1845.909
$ws.Range("J83").Value = 1811.2222
$ws.Range("L83").Value = 16300.9998
$ws.Range("N83").Value = -25660.9998
$ws.Range("H122").Value = 1479.3143
$ws.Range("I122").Value = 1098.1111
$ws.Range("J122").Value = 1882.9412
$ws.Range("K122").Value = 9882.999900000001
$ws.Range("L122").Value = 16946.4708
$ws.Range("M122").Value = -7432.999900000001
$ws.Range("N122").Value = -21846.4708
$ws.Range("H135").Value = 1517.3334
$ws.Range("I135").Value = 768
$ws.Range("J135").Value = 2266.6667
$ws.Range("K135").Value = 6912
$ws.Range("L135").Value = 20400.0003
$ws.Range("M135").Value = -4377
$ws.Range("N135").Value = -25470.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2596.7778
$ws.Range("I107").Value = 2623.875
$ws.Range("J107").Value = 2380
$ws.Range("K107").Value = 2623.875
$ws.Range("L107").Value = 2380
$ws.Range("M107").Value = -703.875
$ws.Range("N107").Value = -6220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5872.4165
$ws.Range("I7").Value = 7468.5
$ws.Range("J7").Value = 5074.375
$ws.Range("K7").Value = 7468.5
$ws.Range("L7").Value = 5074.375
$ws.Range("M7").Value = -7356.5
$ws.Range("N7").Value = -5298.375
$ws.Range("H68").Value = 1823.2858
$ws.Range("I68").Value = 2178.4443
$ws.Range("J68").Value = 1184
$ws.Range("K68").Value = 2178.4443
$ws.Range("L68").Value = 1184
$ws.Range("M68").Value = -1429.4443
$ws.Range("N68").Value = -2682
$ws.Range("H71").Value = 1823.2858
$ws.Range("I71").Value = 2178.4443
$ws.Range("J71").Value = 1184
$ws.Range("K71").Value = 10892.2215
$ws.Range("L71").Value = 5920
$ws.Range("M71").Value = -7148.2215
$ws.Range("N71").Value = -13408
$ws.Range("H100").Value = 2278
$ws.Range("I100").Value = 2286.1428
$ws.Range("J100").Value = 2249.5
$ws.Range("K100").Value = 2286.1428
$ws.Range("L100").Value = 2249.5
$ws.Range("M100").Value = -1745.1428
$ws.Range("N100").Value = -3331.5
$ws.Range("H126").Value = 5872.4165
$ws.Range("I126").Value = 7468.5
$ws.Range("J126").Value = 5074.375
$ws.Range("K126").Value = 22405.5
$ws.Range("L126").Value = 15223.125
$ws.Range("M126").Value = -19935.5
$ws.Range("N126").Value = -20163.125
